$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("openbis-data")

# Rename strain values (shared strings get new text, "JJS-" prefix added)
for ($r = 1; $r -le 6; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $v = $cell.Value()
    if ($v -eq "MGP1") { $cell.Value = "JJS-MGP1" }
    elseif ($v -eq "MGP100") { $cell.Value = "JJS-MGP100" }
    elseif ($v -eq "MGP20") { $cell.Value = "JJS-MGP20" }
    elseif ($v -eq "MGP999") { $cell.Value = "JJS-MGP999" }
}

# Collect row 2 data values (C2:U2) to duplicate into new rows
$values = @()
for ($c = 3; $c -le 21; $c++) {
    $values += , $ws.Cells.Item(2, $c).Value()
}

# Add new row 7: MS
$ws.Cells.Item(7, 1).Value = "MS"
$ws.Cells.Item(7, 2).Value = "OD600"
for ($c = 3; $c -le 21; $c++) {
    $ws.Cells.Item(7, $c).Value = $values[$c - 3]
}

# Add new row 8: WT 168 trp+
$ws.Cells.Item(8, 1).Value = "WT 168 trp+"
$ws.Cells.Item(8, 2).Value = "OD600"
for ($c = 3; $c -le 21; $c++) {
    $ws.Cells.Item(8, $c).Value = $values[$c - 3]
}

$ws.Range("A12").Select()
